# Students data export sheet gets three new trailing columns (Payment Mode,
# Admission Amount, Payment Expected Date), the data row grows from one
# sample student to three, and a couple of existing sample values change.
# Every cell in this sheet is stored as text (even the numeric-looking
# registration numbers / phone numbers / amounts), so every value below is
# written with a leading apostrophe to force text entry and then has its
# format cleared so no numeric/date style sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

function Set-TextCell {
    param($Sheet, [string]$Addr, [string]$Val)
    $Sheet.Range($Addr).Value = "'" + $Val
    $Sheet.Range($Addr).ClearFormats()
}

# New column widths for the three added columns (N, O, P). Excel's COM
# ColumnWidth property is in characters and is offset from the OOXML
# <col width> value by 6/7, empirically.
$ws.Columns("N:O").ColumnWidth = 15 - (6/7)
$ws.Columns("P").ColumnWidth = 20 - (6/7)

# ---- Header row (row 1) ----
$headers = @{
    "A1" = "Registration Number"
    "B1" = "Admission Date"
    "C1" = "Student Name"
    "D1" = "Father's Name"
    "E1" = "Address"
    "F1" = "Contact Number"
    "G1" = "Time Slots"
    "H1" = "Shift"
    "I1" = "Seat Number"
    "J1" = "Amount Paid"
    "K1" = "Amount Due"
    "L1" = "Locker Number"
    "M1" = "Fees Paid Till Date"
    "N1" = "Payment Mode"
    "O1" = "Admission Amount"
    "P1" = "Payment Expected Date"
}
foreach ($addr in $headers.Keys) {
    Set-TextCell $ws $addr $headers[$addr]
}

# ---- Row 2: satyam ----
$row2 = @{
    "A2" = "4"
    "B2" = "2025-01-01"
    "C2" = "satyam"
    "D2" = "satyam ke papa"
    "E2" = "ramkrishnanagar"
    "F2" = "7250585057"
    "G2" = "22:00-06:00"
    "H2" = "2"
    "I2" = "4"
    "J2" = "350.00"
    "K2" = "10.00"
    "L2" = "4"
    "M2" = "2025-02-02"
    "N2" = "online"
    "O2" = "350.00"
    "P2" = "2025-03-23"
}
foreach ($addr in $row2.Keys) {
    Set-TextCell $ws $addr $row2[$addr]
}

# ---- Row 3: test ----
$row3 = @{
    "A3" = "1"
    "B3" = "2025-01-02"
    "C3" = "test"
    "D3" = "testfather"
    "E3" = "ramkrishnanagar"
    "F3" = "7250585058"
    "G3" = "06:00-10:00"
    "H3" = "2"
    "I3" = "2"
    "J3" = "350.00"
    "K3" = "11.00"
    "L3" = "2"
    "M3" = "2025-02-02"
    "N3" = "online"
    "O3" = "0.00"
    "P3" = "2025-03-22"
}
foreach ($addr in $row3.Keys) {
    Set-TextCell $ws $addr $row3[$addr]
}

# ---- Row 4: testAgain (no Payment Expected Date) ----
$row4 = @{
    "A4" = "2"
    "B4" = "2025-03-01"
    "C4" = "testAgain"
    "D4" = "testfather"
    "E4" = "ramkrishnanagar"
    "F4" = "7250585051"
    "G4" = "06:00-10:00"
    "H4" = "1"
    "I4" = "0"
    "J4" = "350.00"
    "K4" = "150.00"
    "L4" = "5"
    "M4" = "2025-03-27"
    "N4" = "online"
    "O4" = "350.00"
}
foreach ($addr in $row4.Keys) {
    Set-TextCell $ws $addr $row4[$addr]
}

Write-Output "students_data sheet updated: headers + 3 rows written"
